$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2439.76
$ws.Range("I19").Value = 3331.8823
$ws.Range("J19").Value = 544
$ws.Range("K19").Value = 3331.8823
$ws.Range("L19").Value = 544
$ws.Range("M19").Value = -3156.8823
$ws.Range("N19").Value = -894

$ws.Range("H96").Value = 1445.2778
$ws.Range("I96").Value = 898.9231
$ws.Range("J96").Value = 2865.8
$ws.Range("K96").Value = 2696.7693
$ws.Range("L96").Value = 8597.400000000001
$ws.Range("M96").Value = -1323.7693
$ws.Range("N96").Value = -11343.4

$ws.Range("H116").Value = 38403.242
$ws.Range("I116").Value = 51851.81
$ws.Range("J116").Value = 3100.75
$ws.Range("K116").Value = 51851.81
$ws.Range("L116").Value = 3100.75
$ws.Range("M116").Value = -48409.81
$ws.Range("N116").Value = -9984.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 715.9655
$ws.Range("I2").Value = 674.52
$ws.Range("J2").Value = 975
$ws.Range("K2").Value = 674.52
$ws.Range("L2").Value = 975
$ws.Range("M2").Value = -561.52
$ws.Range("N2").Value = -1201

$ws.Range("H45").Value = 1052.2858
$ws.Range("I45").Value = 873.2
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 873.2
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -496.2
$ws.Range("N45").Value = -2254

$ws.Range("H113").Value = 34897.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 34897.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 34897.5
$ws.Range("N113").Value = -43575.5

$ws.Range("H116").Value = 715.9655
$ws.Range("I116").Value = 674.52
$ws.Range("J116").Value = 975
$ws.Range("K116").Value = 674.52
$ws.Range("L116").Value = 975
$ws.Range("M116").Value = 1619.48
$ws.Range("N116").Value = -5563

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 715.9655
$ws.Range("I3").Value = 674.52
$ws.Range("J3").Value = 975
$ws.Range("K3").Value = 674.52
$ws.Range("L3").Value = 975
$ws.Range("M3").Value = -560.52
$ws.Range("N3").Value = -1203

$ws.Range("H134").Value = 23423.902
$ws.Range("I134").Value = 28764.184
$ws.Range("J134").Value = 7813.846
$ws.Range("K134").Value = 86292.552
$ws.Range("L134").Value = 23441.538
$ws.Range("M134").Value = -83757.552
$ws.Range("N134").Value = -28511.538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 866.3077
$ws.Range("I107").Value = 308.2
$ws.Range("J107").Value = 1215.125
$ws.Range("K107").Value = 308.2
$ws.Range("L107").Value = 1215.125
$ws.Range("M107").Value = 1611.8
$ws.Range("N107").Value = -5055.125

$ws.Range("H134").Value = 1515.0588
$ws.Range("I134").Value = 898.6061
$ws.Range("J134").Value = 2645.2222
$ws.Range("K134").Value = 2695.8183
$ws.Range("L134").Value = 7935.6666
$ws.Range("M134").Value = -160.8182999999999
$ws.Range("N134").Value = -13005.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 55.846153
$ws.Range("I12").Value = 21.285715
$ws.Range("J12").Value = 96.166664
$ws.Range("K12").Value = 63.857145
$ws.Range("L12").Value = 288.499992
$ws.Range("M12").Value = 109.142855
$ws.Range("N12").Value = -634.499992

$ws.Range("H62").Value = 3244.5
$ws.Range("I62").Value = 998
$ws.Range("J62").Value = 3993.3333
$ws.Range("K62").Value = 2994
$ws.Range("L62").Value = 11979.9999
$ws.Range("M62").Value = -2308
$ws.Range("N62").Value = -13351.9999

$ws.Range("H65").Value = 3244.5
$ws.Range("I65").Value = 998
$ws.Range("J65").Value = 3993.3333
$ws.Range("K65").Value = 8982
$ws.Range("L65").Value = 35939.9997
$ws.Range("M65").Value = -5550
$ws.Range("N65").Value = -42803.9997

$ws.Range("H117").Value = 1026.4166
$ws.Range("I117").Value = 389.5
$ws.Range("J117").Value = 1663.3334
$ws.Range("K117").Value = 1168.5
$ws.Range("L117").Value = 4990.0002
$ws.Range("M117").Value = 2273.5
$ws.Range("N117").Value = -11874.0002

$ws.Range("H121").Value = 1654.0358
$ws.Range("I121").Value = 550
$ws.Range("J121").Value = 1838.0416
$ws.Range("K121").Value = 1650
$ws.Range("L121").Value = 5514.1248
$ws.Range("M121").Value = -340
$ws.Range("N121").Value = -8134.1248

$ws.Range("H122").Value = 718.56665
$ws.Range("I122").Value = 363.6
$ws.Range("J122").Value = 1073.5333
$ws.Range("K122").Value = 3272.4
$ws.Range("L122").Value = 9661.7997
$ws.Range("M122").Value = -822.4000000000001
$ws.Range("N122").Value = -14561.7997

$ws.Range("H131").Value = 740.0303
$ws.Range("I131").Value = 383.45456
$ws.Range("J131").Value = 918.3182
$ws.Range("K131").Value = 1150.36368
$ws.Range("L131").Value = 2754.9546
$ws.Range("M131").Value = 3889.63632
$ws.Range("N131").Value = -12834.9546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1380.7778
$ws.Range("I97").Value = 891.3333
$ws.Range("J97").Value = 1992.5834
$ws.Range("K97").Value = 891.3333
$ws.Range("L97").Value = 1992.5834
$ws.Range("M97").Value = -395.3333
$ws.Range("N97").Value = -2984.5834

$ws.Range("H102").Value = 2366.8262
$ws.Range("I102").Value = 2536.6667
$ws.Range("J102").Value = 1755.4
$ws.Range("K102").Value = 2536.6667
$ws.Range("L102").Value = 1755.4
$ws.Range("M102").Value = -914.6667000000002
$ws.Range("N102").Value = -4999.4

$ws.Range("H132").Value = 3553
$ws.Range("I132").Value = 3660.2964
$ws.Range("J132").Value = 3359.8667
$ws.Range("K132").Value = 10980.8892
$ws.Range("L132").Value = 10079.6001
$ws.Range("M132").Value = -8450.889200000001
$ws.Range("N132").Value = -15139.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280

$ws.Range("H16").Value = 1957.3334
$ws.Range("I16").Value = 2160.1333
$ws.Range("J16").Value = 1450.3334
$ws.Range("K16").Value = 2160.1333
$ws.Range("L16").Value = 1450.3334
$ws.Range("M16").Value = -1990.1333
$ws.Range("N16").Value = -1790.3334

$ws.Range("H21").Value = 5250
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5250
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 5250
$ws.Range("N21").Value = -5598

$ws.Range("H25").Value = 9800
$ws.Range("I25").Value = 9800
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 9800
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9570

$ws.Range("H93").Value = 1899
$ws.Range("I93").Value = 1869.7142
$ws.Range("J93").Value = 1940
$ws.Range("K93").Value = 1869.7142
$ws.Range("L93").Value = 1940
$ws.Range("M93").Value = -621.7141999999999
$ws.Range("N93").Value = -4436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 617
$ws.Range("I107").Value = 525.5
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1576.5
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 343.5
$ws.Range("N107").Value = -6240
